$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 9.473419
$ws.Cells.Item(2, 8).Value = 28.420257
$ws.Cells.Item(2, 9).Value = 0.4285508714990221
$ws.Cells.Item(2, 10).Value = 0.4285508714990221
$ws.Cells.Item(2, 13).Value = 2.302714666666666
$ws.Cells.Item(2, 14).Value = 6.908143999999999
$ws.Cells.Item(2, 15).Value = 0.137852406054037
$ws.Cells.Item(2, 16).Value = 0.137852406054037
$ws.Cells.Item(2, 17).Value = 21.81458087477866
$ws.Cells.Item(2, 18).Value = 196.331227873008
$ws.Cells.Item(2, 19).Value = 0.05907676875269462
$ws.Cells.Item(2, 20).Value = 0.05907676875269461

$ws.Cells.Item(3, 7).Value = 9.473419
$ws.Cells.Item(3, 8).Value = 28.420257
$ws.Cells.Item(3, 9).Value = 0.4285508714990221
$ws.Cells.Item(3, 10).Value = 0.4285508714990221
$ws.Cells.Item(3, 13).Value = 5.106707
$ws.Cells.Item(3, 14).Value = 15.320121
$ws.Cells.Item(3, 15).Value = 0.3057138850737592
$ws.Cells.Item(3, 16).Value = 0.3057138850737592
$ws.Cells.Item(3, 17).Value = 48.377975121233
$ws.Cells.Item(3, 18).Value = 435.401776091097
$ws.Cells.Item(3, 19).Value = 0.1310139518777114
$ws.Cells.Item(3, 20).Value = 0.1310139518777114

$ws.Cells.Item(4, 7).Value = 9.473419
$ws.Cells.Item(4, 8).Value = 28.420257
$ws.Cells.Item(4, 9).Value = 0.4285508714990221
$ws.Cells.Item(4, 10).Value = 0.4285508714990221
$ws.Cells.Item(4, 13).Value = 4.286508666666667
$ws.Cells.Item(4, 14).Value = 12.859526
$ws.Cells.Item(4, 15).Value = 0.2566125720330159
$ws.Cells.Item(4, 16).Value = 0.2566125720330158
$ws.Cells.Item(4, 17).Value = 40.60789264646467
$ws.Cells.Item(4, 18).Value = 365.471033818182
$ws.Cells.Item(4, 19).Value = 0.1099715413823545
$ws.Cells.Item(4, 20).Value = 0.1099715413823545

$ws.Cells.Item(5, 7).Value = 9.473419
$ws.Cells.Item(5, 8).Value = 28.420257
$ws.Cells.Item(5, 9).Value = 0.4285508714990221
$ws.Cells.Item(5, 10).Value = 0.4285508714990221
$ws.Cells.Item(5, 13).Value = 2.214334666666666
$ws.Cells.Item(5, 14).Value = 6.643003999999999
$ws.Cells.Item(5, 15).Value = 0.1325615222882719
$ws.Cells.Item(5, 16).Value = 0.1325615222882719
$ws.Cells.Item(5, 17).Value = 20.97732010355866
$ws.Cells.Item(5, 18).Value = 188.795880932028
$ws.Cells.Item(5, 19).Value = 0.05680935590387596
$ws.Cells.Item(5, 20).Value = 0.05680935590387596

$ws.Cells.Item(6, 7).Value = 9.473419
$ws.Cells.Item(6, 8).Value = 28.420257
$ws.Cells.Item(6, 9).Value = 0.4285508714990221
$ws.Cells.Item(6, 10).Value = 0.4285508714990221
$ws.Cells.Item(6, 13).Value = 2.793938666666667
$ws.Cells.Item(6, 14).Value = 8.381816
$ws.Cells.Item(6, 15).Value = 0.1672596145509162
$ws.Cells.Item(6, 16).Value = 0.1672596145509161
$ws.Cells.Item(6, 17).Value = 26.46815164963467
$ws.Cells.Item(6, 18).Value = 238.213364846712
$ws.Cells.Item(6, 19).Value = 0.07167925358238564
$ws.Cells.Item(6, 20).Value = 0.07167925358238564

$ws.Cells.Item(7, 7).Value = 3.310927333333333
$ws.Cells.Item(7, 8).Value = 9.932782
$ws.Cells.Item(7, 9).Value = 0.1497770545322584
$ws.Cells.Item(7, 10).Value = 0.1497770545322584
$ws.Cells.Item(7, 13).Value = 2.302714666666666
$ws.Cells.Item(7, 14).Value = 6.908143999999999
$ws.Cells.Item(7, 15).Value = 0.137852406054037
$ws.Cells.Item(7, 16).Value = 0.137852406054037
$ws.Cells.Item(7, 17).Value = 7.624120930734221
$ws.Cells.Item(7, 18).Value = 68.617088376608
$ws.Cells.Item(7, 19).Value = 0.02064712733895853
$ws.Cells.Item(7, 20).Value = 0.02064712733895853

$ws.Cells.Item(8, 7).Value = 3.310927333333333
$ws.Cells.Item(8, 8).Value = 9.932782
$ws.Cells.Item(8, 9).Value = 0.1497770545322584
$ws.Cells.Item(8, 10).Value = 0.1497770545322584
$ws.Cells.Item(8, 13).Value = 5.106707
$ws.Cells.Item(8, 14).Value = 15.320121
$ws.Cells.Item(8, 15).Value = 0.3057138850737592
$ws.Cells.Item(8, 16).Value = 0.3057138850737592
$ws.Cells.Item(8, 17).Value = 16.90793578962467
$ws.Cells.Item(8, 18).Value = 152.171422106622
$ws.Cells.Item(8, 19).Value = 0.04578892523596102
$ws.Cells.Item(8, 20).Value = 0.04578892523596102

$ws.Cells.Item(9, 7).Value = 3.310927333333333
$ws.Cells.Item(9, 8).Value = 9.932782
$ws.Cells.Item(9, 9).Value = 0.1497770545322584
$ws.Cells.Item(9, 10).Value = 0.1497770545322584
$ws.Cells.Item(9, 13).Value = 4.286508666666667
$ws.Cells.Item(9, 14).Value = 12.859526
$ws.Cells.Item(9, 15).Value = 0.2566125720330159
$ws.Cells.Item(9, 16).Value = 0.2566125720330158
$ws.Cells.Item(9, 17).Value = 14.19231870903689
$ws.Cells.Item(9, 18).Value = 127.730868381332
$ws.Cells.Item(9, 19).Value = 0.03843467519505211
$ws.Cells.Item(9, 20).Value = 0.0384346751950521

$ws.Cells.Item(10, 7).Value = 3.310927333333333
$ws.Cells.Item(10, 8).Value = 9.932782
$ws.Cells.Item(10, 9).Value = 0.1497770545322584
$ws.Cells.Item(10, 10).Value = 0.1497770545322584
$ws.Cells.Item(10, 13).Value = 2.214334666666666
$ws.Cells.Item(10, 14).Value = 6.643003999999999
$ws.Cells.Item(10, 15).Value = 0.1325615222882719
$ws.Cells.Item(10, 16).Value = 0.1325615222882719
$ws.Cells.Item(10, 17).Value = 7.331501173014221
$ws.Cells.Item(10, 18).Value = 65.98351055712799
$ws.Cells.Item(10, 19).Value = 0.01985467435264969
$ws.Cells.Item(10, 20).Value = 0.01985467435264969

$ws.Cells.Item(11, 7).Value = 3.310927333333333
$ws.Cells.Item(11, 8).Value = 9.932782
$ws.Cells.Item(11, 9).Value = 0.1497770545322584
$ws.Cells.Item(11, 10).Value = 0.1497770545322584
$ws.Cells.Item(11, 13).Value = 2.793938666666667
$ws.Cells.Item(11, 14).Value = 8.381816
$ws.Cells.Item(11, 15).Value = 0.1672596145509162
$ws.Cells.Item(11, 16).Value = 0.1672596145509161
$ws.Cells.Item(11, 17).Value = 9.250527899123556
$ws.Cells.Item(11, 18).Value = 83.254751092112
$ws.Cells.Item(11, 19).Value = 0.0250516524096371
$ws.Cells.Item(11, 20).Value = 0.02505165240963709

$ws.Cells.Item(12, 7).Value = 5.214334666666667
$ws.Cells.Item(12, 8).Value = 15.643004
$ws.Cells.Item(12, 9).Value = 0.2358818569818946
$ws.Cells.Item(12, 10).Value = 0.2358818569818946
$ws.Cells.Item(12, 13).Value = 2.302714666666666
$ws.Cells.Item(12, 14).Value = 6.908143999999999
$ws.Cells.Item(12, 15).Value = 0.137852406054037
$ws.Cells.Item(12, 16).Value = 0.137852406054037
$ws.Cells.Item(12, 17).Value = 12.00712491384178
$ws.Cells.Item(12, 18).Value = 108.064124224576
$ws.Cells.Item(12, 19).Value = 0.03251688152944841
$ws.Cells.Item(12, 20).Value = 0.03251688152944841

$ws.Cells.Item(13, 7).Value = 5.214334666666667
$ws.Cells.Item(13, 8).Value = 15.643004
$ws.Cells.Item(13, 9).Value = 0.2358818569818946
$ws.Cells.Item(13, 10).Value = 0.2358818569818946
$ws.Cells.Item(13, 13).Value = 5.106707
$ws.Cells.Item(13, 14).Value = 15.320121
$ws.Cells.Item(13, 15).Value = 0.3057138850737592
$ws.Cells.Item(13, 16).Value = 0.3057138850737592
$ws.Cells.Item(13, 17).Value = 26.62807934260934
$ws.Cells.Item(13, 18).Value = 239.652714083484
$ws.Cells.Item(13, 19).Value = 0.07211235891634783
$ws.Cells.Item(13, 20).Value = 0.07211235891634782

$ws.Cells.Item(14, 7).Value = 5.214334666666667
$ws.Cells.Item(14, 8).Value = 15.643004
$ws.Cells.Item(14, 9).Value = 0.2358818569818946
$ws.Cells.Item(14, 10).Value = 0.2358818569818946
$ws.Cells.Item(14, 13).Value = 4.286508666666667
$ws.Cells.Item(14, 14).Value = 12.859526
$ws.Cells.Item(14, 15).Value = 0.2566125720330159
$ws.Cells.Item(14, 16).Value = 0.2566125720330158
$ws.Cells.Item(14, 17).Value = 22.35129073956712
$ws.Cells.Item(14, 18).Value = 201.161616656104
$ws.Cells.Item(14, 19).Value = 0.06053025001604797
$ws.Cells.Item(14, 20).Value = 0.06053025001604796

$ws.Cells.Item(15, 7).Value = 5.214334666666667
$ws.Cells.Item(15, 8).Value = 15.643004
$ws.Cells.Item(15, 9).Value = 0.2358818569818946
$ws.Cells.Item(15, 10).Value = 0.2358818569818946
$ws.Cells.Item(15, 13).Value = 2.214334666666666
$ws.Cells.Item(15, 14).Value = 6.643003999999999
$ws.Cells.Item(15, 15).Value = 0.1325615222882719
$ws.Cells.Item(15, 16).Value = 0.1325615222882719
$ws.Cells.Item(15, 17).Value = 11.54628201600178
$ws.Cells.Item(15, 18).Value = 103.916538144016
$ws.Cells.Item(15, 19).Value = 0.03126885804170439
$ws.Cells.Item(15, 20).Value = 0.03126885804170439

$ws.Cells.Item(16, 7).Value = 5.214334666666667
$ws.Cells.Item(16, 8).Value = 15.643004
$ws.Cells.Item(16, 9).Value = 0.2358818569818946
$ws.Cells.Item(16, 10).Value = 0.2358818569818946
$ws.Cells.Item(16, 13).Value = 2.793938666666667
$ws.Cells.Item(16, 14).Value = 8.381816
$ws.Cells.Item(16, 15).Value = 0.1672596145509162
$ws.Cells.Item(16, 16).Value = 0.1672596145509161
$ws.Cells.Item(16, 17).Value = 14.56853124614045
$ws.Cells.Item(16, 18).Value = 131.116781215264
$ws.Cells.Item(16, 19).Value = 0.03945350847834603
$ws.Cells.Item(16, 20).Value = 0.03945350847834602

$ws.Cells.Item(17, 7).Value = 0.4724623333333334
$ws.Cells.Item(17, 8).Value = 1.417387
$ws.Cells.Item(17, 9).Value = 0.02137286915109123
$ws.Cells.Item(17, 10).Value = 0.02137286915109123
$ws.Cells.Item(17, 13).Value = 2.302714666666666
$ws.Cells.Item(17, 14).Value = 6.908143999999999
$ws.Cells.Item(17, 15).Value = 0.137852406054037
$ws.Cells.Item(17, 16).Value = 0.137852406054037
$ws.Cells.Item(17, 17).Value = 1.087945944414222
$ws.Cells.Item(17, 18).Value = 9.791513499728
$ws.Cells.Item(17, 19).Value = 0.002946301436756028
$ws.Cells.Item(17, 20).Value = 0.002946301436756028

$ws.Cells.Item(18, 7).Value = 0.4724623333333334
$ws.Cells.Item(18, 8).Value = 1.417387
$ws.Cells.Item(18, 9).Value = 0.02137286915109123
$ws.Cells.Item(18, 10).Value = 0.02137286915109123
$ws.Cells.Item(18, 13).Value = 5.106707
$ws.Cells.Item(18, 14).Value = 15.320121
$ws.Cells.Item(18, 15).Value = 0.3057138850737592
$ws.Cells.Item(18, 16).Value = 0.3057138850737592
$ws.Cells.Item(18, 17).Value = 2.412726704869667
$ws.Cells.Item(18, 18).Value = 21.714540343827
$ws.Cells.Item(18, 19).Value = 0.006533982863353197
$ws.Cells.Item(18, 20).Value = 0.006533982863353196

$ws.Cells.Item(19, 7).Value = 0.4724623333333334
$ws.Cells.Item(19, 8).Value = 1.417387
$ws.Cells.Item(19, 9).Value = 0.02137286915109123
$ws.Cells.Item(19, 10).Value = 0.02137286915109123
$ws.Cells.Item(19, 13).Value = 4.286508666666667
$ws.Cells.Item(19, 14).Value = 12.859526
$ws.Cells.Item(19, 15).Value = 0.2566125720330159
$ws.Cells.Item(19, 16).Value = 0.2566125720330158
$ws.Cells.Item(19, 17).Value = 2.025213886506889
$ws.Cells.Item(19, 18).Value = 18.226924978562
$ws.Cells.Item(19, 19).Value = 0.00548454692458662
$ws.Cells.Item(19, 20).Value = 0.005484546924586618

$ws.Cells.Item(20, 7).Value = 0.4724623333333334
$ws.Cells.Item(20, 8).Value = 1.417387
$ws.Cells.Item(20, 9).Value = 0.02137286915109123
$ws.Cells.Item(20, 10).Value = 0.02137286915109123
$ws.Cells.Item(20, 13).Value = 2.214334666666666
$ws.Cells.Item(20, 14).Value = 6.643003999999999
$ws.Cells.Item(20, 15).Value = 0.1325615222882719
$ws.Cells.Item(20, 16).Value = 0.1325615222882719
$ws.Cells.Item(20, 17).Value = 1.046189723394222
$ws.Cells.Item(20, 18).Value = 9.415707510548
$ws.Cells.Item(20, 19).Value = 0.002833220070336698
$ws.Cells.Item(20, 20).Value = 0.002833220070336698

$ws.Cells.Item(21, 7).Value = 0.4724623333333334
$ws.Cells.Item(21, 8).Value = 1.417387
$ws.Cells.Item(21, 9).Value = 0.02137286915109123
$ws.Cells.Item(21, 10).Value = 0.02137286915109123
$ws.Cells.Item(21, 13).Value = 2.793938666666667
$ws.Cells.Item(21, 14).Value = 8.381816
$ws.Cells.Item(21, 15).Value = 0.1672596145509162
$ws.Cells.Item(21, 16).Value = 0.1672596145509161
$ws.Cells.Item(21, 17).Value = 1.320030781643556
$ws.Cells.Item(21, 18).Value = 11.880277034792
$ws.Cells.Item(21, 19).Value = 0.003574817856058685
$ws.Cells.Item(21, 20).Value = 0.003574817856058685

$ws.Cells.Item(22, 7).Value = 3.634561333333334
$ws.Cells.Item(22, 8).Value = 10.903684
$ws.Cells.Item(22, 9).Value = 0.1644173478357336
$ws.Cells.Item(22, 10).Value = 0.1644173478357336
$ws.Cells.Item(22, 13).Value = 2.302714666666666
$ws.Cells.Item(22, 14).Value = 6.908143999999999
$ws.Cells.Item(22, 15).Value = 0.137852406054037
$ws.Cells.Item(22, 16).Value = 0.137852406054037
$ws.Cells.Item(22, 17).Value = 8.369357689166222
$ws.Cells.Item(22, 18).Value = 75.32421920249601
$ws.Cells.Item(22, 19).Value = 0.02266532699617939
$ws.Cells.Item(22, 20).Value = 0.02266532699617939

$ws.Cells.Item(23, 7).Value = 3.634561333333334
$ws.Cells.Item(23, 8).Value = 10.903684
$ws.Cells.Item(23, 9).Value = 0.1644173478357336
$ws.Cells.Item(23, 10).Value = 0.1644173478357336
$ws.Cells.Item(23, 13).Value = 5.106707
$ws.Cells.Item(23, 14).Value = 15.320121
$ws.Cells.Item(23, 15).Value = 0.3057138850737592
$ws.Cells.Item(23, 16).Value = 0.3057138850737592
$ws.Cells.Item(23, 17).Value = 18.56063980286267
$ws.Cells.Item(23, 18).Value = 167.045758225764
$ws.Cells.Item(23, 19).Value = 0.05026466618038577
$ws.Cells.Item(23, 20).Value = 0.05026466618038576

$ws.Cells.Item(24, 7).Value = 3.634561333333334
$ws.Cells.Item(24, 8).Value = 10.903684
$ws.Cells.Item(24, 9).Value = 0.1644173478357336
$ws.Cells.Item(24, 10).Value = 0.1644173478357336
$ws.Cells.Item(24, 13).Value = 4.286508666666667
$ws.Cells.Item(24, 14).Value = 12.859526
$ws.Cells.Item(24, 15).Value = 0.2566125720330159
$ws.Cells.Item(24, 16).Value = 0.2566125720330158
$ws.Cells.Item(24, 17).Value = 15.57957865486489
$ws.Cells.Item(24, 18).Value = 140.216207893784
$ws.Cells.Item(24, 19).Value = 0.04219155851497462
$ws.Cells.Item(24, 20).Value = 0.04219155851497462

$ws.Cells.Item(25, 7).Value = 3.634561333333334
$ws.Cells.Item(25, 8).Value = 10.903684
$ws.Cells.Item(25, 9).Value = 0.1644173478357336
$ws.Cells.Item(25, 10).Value = 0.1644173478357336
$ws.Cells.Item(25, 13).Value = 2.214334666666666
$ws.Cells.Item(25, 14).Value = 6.643003999999999
$ws.Cells.Item(25, 15).Value = 0.1325615222882719
$ws.Cells.Item(25, 16).Value = 0.1325615222882719
$ws.Cells.Item(25, 17).Value = 8.048135158526222
$ws.Cells.Item(25, 18).Value = 72.43321642673601
$ws.Cells.Item(25, 19).Value = 0.02179541391970516
$ws.Cells.Item(25, 20).Value = 0.02179541391970516

$ws.Cells.Item(26, 7).Value = 3.634561333333334
$ws.Cells.Item(26, 8).Value = 10.903684
$ws.Cells.Item(26, 9).Value = 0.1644173478357336
$ws.Cells.Item(26, 10).Value = 0.1644173478357336
$ws.Cells.Item(26, 13).Value = 2.793938666666667
$ws.Cells.Item(26, 14).Value = 8.381816
$ws.Cells.Item(26, 15).Value = 0.1672596145509162
$ws.Cells.Item(26, 16).Value = 0.1672596145509161
$ws.Cells.Item(26, 17).Value = 10.15474144557156
$ws.Cells.Item(26, 18).Value = 91.39267301014402
$ws.Cells.Item(26, 19).Value = 0.02750038222448872
$ws.Cells.Item(26, 20).Value = 0.02750038222448872
